# Fix Training Data Issue (#48)
#
# The "Date" column (BF) on the sheet held the literal text "5-28-2011-12"
# for every team row. Because of how the NBA stats site displayed dates,
# the data was actually for 2012-05-28, one day off from what the label
# implied. Correct every BF2:BF31 cell to read "2012-05-28" (still as
# literal text, not a real Excel date/number) so downstream model-training
# code keeps treating the column as a plain string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "5-28-2011-12"
$newDate = "2012-05-28"

# Excel "helpfully" reinterprets a bare assignment of an ISO-looking string
# (e.g. "2012-05-28") as a real date serial number, which would change both
# the stored type and the cell style. To keep the replacement cells as
# literal text (same as every other cell in this column), stage the
# corrected string in a scratch cell using a leading apostrophe (forces
# text entry), then copy/paste-values it into each cell that still holds
# the old label.
$helper = $ws.Range("A100")
$helper.Value = "'" + $newDate
$helper.Copy()

$lastRow = 31
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Range("BF" + $row)
    if ($cell.Value() -eq $oldDate) {
        $cell.PasteSpecial(-4163)  # xlPasteValues - values only, no formatting/style carried over
    }
}

# Clean up the scratch row so the sheet's used range / dimension is
# unaffected by the helper cell.
$helper.EntireRow.Delete()
$excel.CutCopyMode = $false
